# Update recalculated TPM-derived statistics for the Ncam1-Robo1 LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05354133333333333
$ws.Range("H2").Value = 0.160624
$ws.Range("I2").Value = 0.00209946492164722
$ws.Range("J2").Value = 0.00209946492164722
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.007258333333333333
$ws.Range("N2").Value = 0.021775
$ws.Range("O2").Value = 0.000328667160253549
$ws.Range("P2").Value = 0.000328667160253549
$ws.Range("Q2").Value = 0.0003886208444444444
$ws.Range("R2").Value = 0.0034975876
$ws.Range("S2").Value = 0.000000690025173849731493632906
$ws.Range("T2").Value = 0.000000690025173849731387753787

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05354133333333333
$ws.Range("H3").Value = 0.160624
$ws.Range("I3").Value = 0.00209946492164722
$ws.Range("J3").Value = 0.00209946492164722
$ws.Range("O3").Value = 0.7778551418094273
$ws.Range("P3").Value = 0.7778551418094272
$ws.Range("Q3").Value = 0.9197472659946666
$ws.Range("R3").Value = 8.277725393952
$ws.Range("S3").Value = 0.001633079584351816
$ws.Range("T3").Value = 0.001633079584351816

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05354133333333333
$ws.Range("H4").Value = 0.160624
$ws.Range("I4").Value = 0.00209946492164722
$ws.Range("J4").Value = 0.00209946492164722
$ws.Range("M4").Value = 4.898620999999999
$ws.Range("N4").Value = 14.695863
$ws.Range("O4").Value = 0.2218161910303192
$ws.Range("P4").Value = 0.2218161910303192
$ws.Range("Q4").Value = 0.2622786998346666
$ws.Range("R4").Value = 2.360508298512
$ws.Range("S4").Value = 0.0004656953121215539
$ws.Range("T4").Value = 0.0004656953121215538

# Row 5
$ws.Range("I5").Value = 0.05460670042535784
$ws.Range("J5").Value = 0.05460670042535784
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.007258333333333333
$ws.Range("N5").Value = 0.021775
$ws.Range("O5").Value = 0.000328667160253549
$ws.Range("P5").Value = 0.000328667160253549
$ws.Range("Q5").Value = 0.01010795741944445
$ws.Range("R5").Value = 0.09097161677500001
$ws.Range("S5").Value = 0.000017947429159618628712177454
$ws.Range("T5").Value = 0.000017947429159618628712177454

# Row 6
$ws.Range("I6").Value = 0.05460670042535784
$ws.Range("J6").Value = 0.05460670042535784
$ws.Range("O6").Value = 0.7778551418094273
$ws.Range("P6").Value = 0.7778551418094272
$ws.Range("S6").Value = 0.04247610270311163
$ws.Range("T6").Value = 0.04247610270311163

# Row 7
$ws.Range("I7").Value = 0.05460670042535784
$ws.Range("J7").Value = 0.05460670042535784
$ws.Range("M7").Value = 4.898620999999999
$ws.Range("N7").Value = 14.695863
$ws.Range("O7").Value = 0.2218161910303192
$ws.Range("P7").Value = 0.2218161910303192
$ws.Range("Q7").Value = 6.821821237473666
$ws.Range("R7").Value = 61.396391137263
$ws.Range("S7").Value = 0.01211265029308659
$ws.Range("T7").Value = 0.01211265029308659

# Row 8
$ws.Range("G8").Value = 24.05622933333333
$ws.Range("H8").Value = 72.168688
$ws.Range("I8").Value = 0.943293834652995
$ws.Range("J8").Value = 0.943293834652995
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.007258333333333333
$ws.Range("N8").Value = 0.021775
$ws.Range("O8").Value = 0.000328667160253549
$ws.Range("P8").Value = 0.000328667160253549
$ws.Range("Q8").Value = 0.1746081312444444
$ws.Range("R8").Value = 1.5714731812
$ws.Range("S8").Value = 0.0003100297059200807
$ws.Range("T8").Value = 0.0003100297059200806

# Row 9
$ws.Range("G9").Value = 24.05622933333333
$ws.Range("H9").Value = 72.168688
$ws.Range("I9").Value = 0.943293834652995
$ws.Range("J9").Value = 0.943293834652995
$ws.Range("O9").Value = 0.7778551418094273
$ws.Range("P9").Value = 0.7778551418094272
$ws.Range("Q9").Value = 413.2443064450027
$ws.Range("R9").Value = 3719.198758005025
$ws.Range("S9").Value = 0.7337459595219639
$ws.Range("T9").Value = 0.7337459595219638

# Row 10
$ws.Range("G10").Value = 24.05622933333333
$ws.Range("H10").Value = 72.168688
$ws.Range("I10").Value = 0.943293834652995
$ws.Range("J10").Value = 0.943293834652995
$ws.Range("M10").Value = 4.898620999999999
$ws.Range("N10").Value = 14.695863
$ws.Range("O10").Value = 0.2218161910303192
$ws.Range("P10").Value = 0.2218161910303192
$ws.Range("Q10").Value = 117.8423501930827
$ws.Range("R10").Value = 1060.581151737744
$ws.Range("S10").Value = 0.2092378454251111
$ws.Range("T10").Value = 0.2092378454251111
